$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.8377025
$ws.Range("N2").Value = 5.675405
$ws.Range("O2").Value = 0.08520139853031897
$ws.Range("P2").Value = 0.05979321158534227
$ws.Range("Q2").Value = 1.826275332338333
$ws.Range("R2").Value = 10.95765199403
$ws.Range("S2").Value = 0.08520139853031897
$ws.Range("T2").Value = 0.05979321158534227

# Row 3
$ws.Range("O3").Value = 0.1519846751111432
$ws.Range("P3").Value = 0.1599912441592654
$ws.Range("S3").Value = 0.1519846751111432
$ws.Range("T3").Value = 0.1599912441592654

# Row 4
$ws.Range("M4").Value = 4.372280666666666
$ws.Range("N4").Value = 13.116842
$ws.Range("O4").Value = 0.1312767732230829
$ws.Range("P4").Value = 0.1381924477702479
$ws.Range("Q4").Value = 2.813891987476888
$ws.Range("R4").Value = 25.32502788729199
$ws.Range("S4").Value = 0.1312767732230829
$ws.Range("T4").Value = 0.1381924477702479

# Row 5
$ws.Range("M5").Value = 2.1625465
$ws.Range("N5").Value = 4.325093
$ws.Range("O5").Value = 0.06492998691259792
$ws.Range("P5").Value = 0.04556700374251401
$ws.Range("Q5").Value = 1.391761584586333
$ws.Range("R5").Value = 8.350569507517999
$ws.Range("S5").Value = 0.06492998691259792
$ws.Range("T5").Value = 0.04556700374251401

# Row 6
$ws.Range("M6").Value = 15.97657333333333
$ws.Range("N6").Value = 47.92972
$ws.Range("O6").Value = 0.479693129114909
$ws.Range("P6").Value = 0.5049634148023289
$ws.Range("Q6").Value = 10.28212850852444
$ws.Range("R6").Value = 92.53915657672
$ws.Range("S6").Value = 0.479693129114909
$ws.Range("T6").Value = 0.5049634148023289

# Row 7
$ws.Range("M7").Value = 2.894743333333333
$ws.Range("N7").Value = 8.684229999999999
$ws.Range("O7").Value = 0.08691403710794816
$ws.Range("P7").Value = 0.0914926779403015
$ws.Range("Q7").Value = 1.862985405664444
$ws.Range("R7").Value = 16.76686865098
$ws.Range("S7").Value = 0.08691403710794816
$ws.Range("T7").Value = 0.0914926779403015

$wb.Save()
